# "Generate Report for Handoff" - refresh the localization-status report.
#
# For the rows that were still at "Ready for handoff" (the four newest
# entries, rows 4-7) on each language sheet, the handoff run:
#   - bumps the Priority from "low" to "ht" (handed off), and
#   - stamps the Latest Handoff Datetime with the new xliff generation time.

$wb = $excel.ActiveWorkbook

# zh-cn: handoff xliffs generated at 2016-08-17 22:30:48
$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($row in 4..7) {
    $wsZhCn.Range("E$row").Value = "ht"
    $wsZhCn.Range("H$row").Value = "2016-08-17 22:30:48"
}

# de-de: handoff xliffs generated at 2016-08-17 22:30:54
$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($row in 4..7) {
    $wsDeDe.Range("E$row").Value = "ht"
    $wsDeDe.Range("H$row").Value = "2016-08-17 22:30:54"
}

# Overview's "Latest HO Xliff Generate Date" mirrors de-de's handoff datetime
# (it is the same underlying reference-language timestamp), so refresh it too.
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($row in 4..7) {
    $wsOverview.Range("G$row").Value = "2016-08-17 22:30:54"
}
